# Update countries & provincias Spain
# Applies updated COVID-19 statistics to the "Pais" sheet and swaps the
# Togo / Suazilandia rows (Suazilandia's totals overtook Togo's, so the
# sheet - sorted descending by "Casos totales" - now lists Suazilandia
# before Togo).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Estados Unidos (row 4) ---
$ws.Range("B4").Value = 1372740
$ws.Range("C4").Value = 5102
$ws.Range("D4").Value = 256972
$ws.Range("E4").Value = 1034835
$ws.Range("F4").Value = 16514
$ws.Range("G4").Value = 146
$ws.Range("H4").Value = 80933

# --- Italia (row 8) ---
$ws.Range("B8").Value = 219814
$ws.Range("C8").Value = 744
$ws.Range("D8").Value = 106587
$ws.Range("E8").Value = 82488
$ws.Range("F8").Value = 999
$ws.Range("G8").Value = 179
$ws.Range("H8").Value = 30739

# --- India (row 15) ---
$ws.Range("B15").Value = 69400
$ws.Range("C15").Value = 2239
$ws.Range("D15").Value = 21664
$ws.Range("E15").Value = 45482
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 42
$ws.Range("H15").Value = 2254

# --- Canada (row 16) ---
$ws.Range("B16").Value = 69157
$ws.Range("C16").Value = 309
$ws.Range("D16").Value = 32096
$ws.Range("E16").Value = 32154
$ws.Range("F16").Value = 502
$ws.Range("G16").Value = 37
$ws.Range("H16").Value = 4907

# --- Suiza (row 23) ---
$ws.Range("B23").Value = 30344
$ws.Range("C23").Value = 39
$ws.Range("D23").Value = 26600
$ws.Range("E23").Value = 1899
$ws.Range("F23").Value = 89
$ws.Range("G23").Value = 12
$ws.Range("H23").Value = 1845

# --- Polonia (row 34) ---
$ws.Range("B34").Value = 16326
$ws.Range("C34").Value = 330
$ws.Range("D34").Value = 5816
$ws.Range("E34").Value = 9699
$ws.Range("F34").Value = 160
$ws.Range("G34").Value = 11
$ws.Range("H34").Value = 811

# --- Rumania (row 39) ---
$ws.Range("B39").Value = 15588
$ws.Range("C39").Value = 226
$ws.Range("D39").Value = 7245
$ws.Range("E39").Value = 7361
$ws.Range("F39").Value = 255
$ws.Range("G39").Value = 21
$ws.Range("H39").Value = 982

# --- Cuba (row 83) ---
$ws.Range("B83").Value = 1783
$ws.Range("C83").Value = 17
$ws.Range("D83").Value = 1229
$ws.Range("E83").Value = 477
$ws.Range("F83").Value = 7
$ws.Range("G83").Value = 0
$ws.Range("H83").Value = 77

# --- Montenegro (row 133) ---
$ws.Range("B133").Value = 324
$ws.Range("C133").Value = 0
$ws.Range("D133").Value = 294
$ws.Range("E133").Value = 21
$ws.Range("F133").Value = 2
$ws.Range("G133").Value = 0
$ws.Range("H133").Value = 9

# --- Suazilandia / Togo swap (rows 149-150) ---
# Row 149 becomes Suazilandia with refreshed totals (now higher than Togo's).
$ws.Range("A149").Value = "Suazilandia"
$ws.Range("B149").Value = 175
$ws.Range("C149").Value = 3
$ws.Range("D149").Value = 28
$ws.Range("E149").Value = 145
$ws.Range("F149").Value = 0
$ws.Range("G149").Value = 0
$ws.Range("H149").Value = 2

# Row 150 becomes Togo, keeping its previous totals.
$ws.Range("A150").Value = "Togo"
$ws.Range("B150").Value = 174
$ws.Range("C150").Value = 0
$ws.Range("D150").Value = 89
$ws.Range("E150").Value = 74
$ws.Range("F150").Value = 0
$ws.Range("G150").Value = 0
$ws.Range("H150").Value = 11

# --- Libia (row 171) ---
$ws.Range("D171").Value = 28
$ws.Range("E171").Value = 33
